# Apply the authored changes to the "range_calculation" workbook (Tabelle1):
#  - Update the resistor/capacitor pair in the last calc block (B35, C35).
#    A35 holds the formula =1/(2*PI()*B35*C35*0.000000001) and recalculates
#    automatically once the inputs change.
#  - Move the active-cell selection from D23 to C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the inputs feeding the A35 formula.
$ws.Range("B35").Value = 100
$ws.Range("C35").Value = 15

# Move the selection/active cell to C23 (was D23).
$ws.Range("C23").Select()
